# The source data added a new daily price record for Mango at
# "Macroferia Regional de Talca" market. In the underlying "logica_diaria"
# (daily logic) consolidation sheet, newer records are inserted at the top
# of the data block (row 116, right after the most recent previously-seen
# entry), pushing the rest of the historical rows down by one.
#
# This script:
#   1. Inserts a new blank row at row 116 (shifting rows 116:188 -> 117:189).
#   2. Populates the new row 116 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at 116, shifting existing rows 116-188 down to 117-189.
$ws.Rows.Item(116).EntireRow.Insert()

# 2. Fill in the new record's data in row 116.
$ws.Range("A116").Value = 5
$ws.Range("B116").Value = "Macroferia Regional de Talca"
$ws.Range("C116").Value = "Maule"
$ws.Range("D116").Value = 45161
$ws.Range("E116").Value = 7
$ws.Range("F116").Value = "Fruta"
$ws.Range("G116").Value = 100108
$ws.Range("H116").Value = "Tropicales y subtropicales"
$ws.Range("I116").Value = 100108002
$ws.Range("J116").Value = "Mango"
$ws.Range("K116").Value = "Sin especificar"
$ws.Range("L116").Value = "Primera"
$ws.Range("M116").Value = 245
$ws.Range("N116").Value = 8000
$ws.Range("O116").Value = 8000
$ws.Range("P116").Value = 8000
$ws.Range("Q116").Value = '$/bandeja 4 kilos'
$ws.Range("R116").Value = "Brasil"
$ws.Range("S116").Value = 2000
$ws.Range("T116").Value = 4
